# "Actualizacion data y funcionalidades"
# The data-driven fixture's "Datos" sheet gets its tipoCuenta value for the
# second test row (N2) updated from "Ahorros" to "Corriente", and the
# workbook's last active selection moves to N8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Update the account-type value used by the second data row.
$ws.Range("N2").Value = "Corriente"

# Leave the cursor where the author left it when they saved the file.
$ws.Range("N8").Select() | Out-Null
